$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to Text format first so the literal digit/dot string is preserved.
$textCells = @("D5", "D6", "D11", "D14", "D15", "D19", "D20", "D22", "D24", "D25", "D27", "D29", "D32", "D34", "D37", "D39", "D41", "D42", "D48", "D49", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "60.665.99"
$ws.Range("E2").Value = "  -5.22%  "
$ws.Range("D3").Value = "3.015.95"
$ws.Range("E3").Value = "  -6.59%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "578.90"
$ws.Range("E5").Value = "  -2.81%  "
$ws.Range("D6").Value = "126.96"
$ws.Range("E6").Value = "  -8.05%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "3.012.21"
$ws.Range("E8").Value = "  -6.71%  "
$ws.Range("E9").Value = "  -3.20%  "
$ws.Range("E10").Value = "  -7.81%  "
$ws.Range("D11").Value = "5.15"
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("E12").Value = "  -3.76%  "
$ws.Range("E13").Value = "  -7.65%  "
$ws.Range("D14").Value = "32.68"
$ws.Range("E14").Value = "  -8.47%  "
$ws.Range("D15").Value = "0.120"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "3.508.82"
$ws.Range("E16").Value = "  -6.82%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "60.642.65"
$ws.Range("E17").Value = "  -5.28%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.003.00"
$ws.Range("E18").Value = "  -7.07%  "
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  -3.60%  "
$ws.Range("D20").Value = "433.17"
$ws.Range("E20").Value = "  -7.65%  "
$ws.Range("E21").Value = "  -7.17%  "
$ws.Range("D22").Value = "0.667"
$ws.Range("E22").Value = "  -6.19%  "
$ws.Range("E23").Value = "  -9.47%  "
$ws.Range("D24").Value = "12.81"
$ws.Range("E24").Value = "  -5.36%  "
$ws.Range("D25").Value = "79.61"
$ws.Range("E25").Value = "  -5.19%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  -4.95%  "
$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  -7.45%  "
$ws.Range("E30").Value = "  -8.44%  "
$ws.Range("E31").Value = "  -10.51%  "
$ws.Range("D32").Value = "25.37"
$ws.Range("E32").Value = "  -8.53%  "
$ws.Range("E33").Value = "  -9.93%  "
$ws.Range("D34").Value = "2.16"
$ws.Range("E34").Value = "  -12.91%  "
$ws.Range("E35").Value = "  -8.40%  "
$ws.Range("E36").Value = "  -5.72%  "
$ws.Range("D37").Value = "50.07"
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("E38").Value = "  -10.77%  "
$ws.Range("D39").Value = "8.43"
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("E40").Value = "  -9.19%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "386.86"
$ws.Range("E42").Value = "  -5.43%  "
$ws.Range("E43").Value = "  -10.40%  "
$ws.Range("D44").Value = "2.661.63"
$ws.Range("E44").Value = "  -7.00%  "
$ws.Range("E46").Value = "  -8.52%  "
$ws.Range("E47").Value = "  -7.61%  "
$ws.Range("D48").Value = "118.67"
$ws.Range("E48").Value = "  -7.84%  "
$ws.Range("D49").Value = "0.108"
$ws.Range("E49").Value = "  -4.46%  "
$ws.Range("E50").Value = "  -8.44%  "
$ws.Range("D51").Value = "0.134"
$ws.Range("E51").Value = "  +2.33%  "

# Restore the default (General) style on the cells we force-formatted,
# now that the literal text is safely stored, so no stray style id lingers.
foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}
